$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.350.66"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").Value = "3.197.23"
$ws.Range("E3").Value = "  +4.98%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'207.92"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").Value = "'632.16"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  +6.57%  "
$ws.Range("D9").Value = "'0.226"
$ws.Range("E9").Value = "  +11.92%  "
$ws.Range("D10").Value = "3.203.04"
$ws.Range("E10").Value = "  +5.23%  "
$ws.Range("D11").Value = "'0.581"
$ws.Range("E11").Value = "  +33.78%  "
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "'5.42"
$ws.Range("E13").Value = "  +7.34%  "
$ws.Range("D14").Value = "3.797.93"
$ws.Range("E14").Value = "  +6.06%  "
$ws.Range("D15").Value = "'0.0000226"
$ws.Range("E15").Value = "  +19.58%  "
$ws.Range("D16").Value = "'31.69"
$ws.Range("E16").Value = "  +8.33%  "
$ws.Range("D17").Value = "79.228.42"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").Value = "3.197.14"
$ws.Range("E18").Value = "  +5.90%  "
$ws.Range("D19").Value = "'14.45"
$ws.Range("E19").Value = "  +5.57%  "
$ws.Range("D20").Value = "'9.46"
$ws.Range("E20").Value = "  +4.39%  "
$ws.Range("D21").Value = "'431.79"
$ws.Range("E21").Value = "  +15.00%  "
$ws.Range("D22").Value = "'2.83"
$ws.Range("E22").Value = "  +24.27%  "
$ws.Range("D23").Value = "'5.01"
$ws.Range("E23").Value = "  +15.84%  "
$ws.Range("E24").Value = "  +5.93%  "
$ws.Range("D25").Value = "3.368.52"
$ws.Range("E25").Value = "  +7.49%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'77.27"
$ws.Range("E26").Value = "  +5.92%  "
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").Value = "'4.76"
$ws.Range("E27").Value = "  +8.48%  "
$ws.Range("D28").Value = "'10.94"
$ws.Range("E28").Value = "  +10.65%  "
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").Value = "'0.0000115"
$ws.Range("E30").Value = "  +5.83%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "'8.97"
$ws.Range("E32").Value = "  +7.62%  "
$ws.Range("D33").Value = "'1.47"
$ws.Range("E33").Value = "  +4.54%  "
$ws.Range("D34").Value = "'516.78"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").Value = "'1.98"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").Value = "'0.129"
$ws.Range("E36").Value = "  +22.90%  "
$ws.Range("D37").Value = "'0.136"
$ws.Range("E37").Value = "  +21.45%  "
$ws.Range("D38").Value = "'22.89"
$ws.Range("E38").Value = "  +9.77%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "'0.408"
$ws.Range("E40").Value = "  +5.32%  "
$ws.Range("D41").Value = "'164.29"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").Value = "'195.89"
$ws.Range("E42").Value = "  +4.40%  "
$ws.Range("D43").Value = "'20.00"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").Value = "'5.44"
$ws.Range("E45").Value = "  +8.34%  "
$ws.Range("D46").Value = "'0.815"
$ws.Range("E46").Value = "  +13.78%  "
$ws.Range("D47").Value = "'1.78"
$ws.Range("E47").Value = "  +7.63%  "
$ws.Range("E48").Value = "  +4.10%  "
$ws.Range("D49").Value = "'43.03"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'2.53"
$ws.Range("E50").Value = "  +9.25%  "
$ws.Range("D51").Value = "'0.629"
$ws.Range("E51").Value = "  +2.60%  "

Write-Host "Applied cryptos update"
